$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "UItests" to "ErrorHandling"
$ws.Name = "ErrorHandling"

# Move/update the selected cell on the sheet
$ws.Range("C18").Select()

$wb.Save()
